# Updates reported symbol list prices / volume labels / swapped two rows
# to match the latest coinranking.com snapshot (Fri Dec 16 15:55:47 UTC 2022).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "247.70"; ForceText = $true }
    @{ Cell = "D3"; Value = "24.11"; ForceText = $true }
    @{ Cell = "D4"; Value = "5.871"; ForceText = $true }
    @{ Cell = "D5"; Value = "0.05884"; ForceText = $true }
    @{ Cell = "D6"; Value = "3.429"; ForceText = $true }
    @{ Cell = "D7"; Value = "6.515"; ForceText = $true }
    @{ Cell = "D8"; Value = "1.329"; ForceText = $true }
    @{ Cell = "D9"; Value = "0.7976"; ForceText = $true }
    @{ Cell = "D10"; Value = "0.1476"; ForceText = $true }
    @{ Cell = "D11"; Value = "0.07754"; ForceText = $true }
    @{ Cell = "D12"; Value = "0.03304"; ForceText = $true }
    @{ Cell = "D13"; Value = "0.03010"; ForceText = $true }
    @{ Cell = "D14"; Value = "0.09221"; ForceText = $true }
    @{ Cell = "D15"; Value = "3.572"; ForceText = $true }
    @{ Cell = "D16"; Value = "0.001664"; ForceText = $true }
    @{ Cell = "D17"; Value = "0.04759"; ForceText = $true }
    @{ Cell = "D18"; Value = "0.0006036"; ForceText = $true }
    @{ Cell = "E18"; Value = "17OneONE"; ForceText = $false }
    @{ Cell = "D19"; Value = "0.006240"; ForceText = $true }
    @{ Cell = "D20"; Value = "0.005533"; ForceText = $true }
    @{ Cell = "D21"; Value = "0.001069"; ForceText = $true }
    @{ Cell = "D22"; Value = "0.0001501"; ForceText = $true }
    @{ Cell = "D23"; Value = "3.698"; ForceText = $true }
    @{ Cell = "D25"; Value = "0.3350"; ForceText = $true }
    @{ Cell = "D26"; Value = "0.1254"; ForceText = $true }
    @{ Cell = "D27"; Value = "0.0006275"; ForceText = $true }
    @{ Cell = "D40"; Value = "0.04364"; ForceText = $true }
    @{ Cell = "D41"; Value = "0.007045"; ForceText = $true }
    @{ Cell = "B42"; Value = "CEJI"; ForceText = $false }
    @{ Cell = "C42"; Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"; ForceText = $false }
    @{ Cell = "D42"; Value = "0.003603"; ForceText = $true }
    @{ Cell = "E42"; Value = "41CEJICEJI"; ForceText = $false }
    @{ Cell = "B43"; Value = "BKEXToken"; ForceText = $false }
    @{ Cell = "C43"; Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"; ForceText = $false }
    @{ Cell = "D43"; Value = "0.1061"; ForceText = $true }
    @{ Cell = "E43"; Value = "42BKEXTokenBKK"; ForceText = $false }
    @{ Cell = "D44"; Value = "0.009650"; ForceText = $true }
    @{ Cell = "E45"; Value = "44ACDXExchangeACXTBestin24h"; ForceText = $false }
    @{ Cell = "D46"; Value = "0.00005891"; ForceText = $true }
    @{ Cell = "D48"; Value = "0.9909"; ForceText = $true }
    @{ Cell = "D49"; Value = "0.1088"; ForceText = $true }
    @{ Cell = "E49"; Value = "48BOLOBOLOWorstin24h"; ForceText = $false }
    @{ Cell = "D50"; Value = "0.00002102"; ForceText = $true }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    if ($u.ForceText) {
        # Preserve the original inline-string representation (e.g. "247.70")
        # instead of letting Excel coerce the literal into a Number.
        $cell.NumberFormat = "@"
        $cell.Value = $u.Value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $u.Value
    }
}
